$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before G (shifts old G "Description" data to H)
$ws.Columns("G").Insert()
$ws.Columns("G").ColumnWidth = $ws.Columns("F").ColumnWidth

# Fill in the new row 10 data for the "assault" command
$ws.Range("B10").Value = "Basic"
$ws.Range("C10").Value = "Combat"
$ws.Range("D10").Value = "Yes"
$ws.Range("E10").Value = "Yes"
$ws.Range("F10").Value = "Yes"
$ws.Range("H10").Value = "From MaritimeUnit to LandUnits"

# Add header for the new column G, matching the header style of A1
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "BTUs"

# Update selection
$ws.Range("G2").Select()
